$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").NumberFormat = "@"
$ws.Range("A24").Value = "28"
$ws.Range("A24").Style = "Normal"

$ws.Range("B24").Value = "[BUG] <title>"
$ws.Range("C24").Value = "open"
$ws.Range("D24").Value = "2025-03-26T07:00:05Z"
$ws.Range("E24").Value = "bug"
